# Fruta / hortaliza, semanal
# Insert a new weekly record before the existing row 18 ("Terminal
# Hortofrutícola Agro Chillán" - Mango sheet). Inserting the row shifts the
# previous rows 18-40 down to 19-41 (carrying all of their data along), and
# the freed-up row 18 is then populated with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 18..40 down to 19..41, leaving a blank row 18 to fill in below.
$ws.Rows.Item(18).Insert()

# New row 18 content (mirrors the surrounding rows' fixed columns).
$ws.Cells.Item(18, 1).Value = 7
$ws.Cells.Item(18, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(18, 3).Value = "Ñuble"
$ws.Cells.Item(18, 4).Value = "8/30/2021"
$ws.Cells.Item(18, 5).Value = 16
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100108
$ws.Cells.Item(18, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(18, 9).Value = 100108002
$ws.Cells.Item(18, 10).Value = "Mango"
$ws.Cells.Item(18, 11).Value = "Sin especificar"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 60
$ws.Cells.Item(18, 14).Value = 8500
$ws.Cells.Item(18, 15).Value = 9000
$ws.Cells.Item(18, 16).Value = 8750
$ws.Cells.Item(18, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(18, 18).Value = "Perú"
$ws.Cells.Item(18, 19).Value = 2188
$ws.Cells.Item(18, 20).Value = 4
